# Natmi LR-pair results (Angptl1-Tek): rebuild the results table to include
# the "ECs" sending/target cluster alongside the existing "FAPs" and "sCs"
# clusters (full 3x3 cross product of sending x target clusters), per
# Dr Hou's advice. Ligand symbol is always "Angptl1" and receptor symbol is
# always "Tek".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Angptl1"
$ws.Range("C2").Value = "Tek"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1275283333333333
$ws.Range("H2").Value = 0.382585
$ws.Range("I2").Value = 0.001938448326189378
$ws.Range("J2").Value = 0.001938448326189378
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 34.10446833333334
$ws.Range("N2").Value = 102.313405
$ws.Range("O2").Value = 0.5118942073015388
$ws.Range("P2").Value = 0.5118942073015389
$ws.Range("Q2").Value = 4.349286005769445
$ws.Range("R2").Value = 39.143574051925
$ws.Range("S2").Value = 0.0009922804693297063
$ws.Range("T2").Value = 0.0009922804693297065

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Angptl1"
$ws.Range("C3").Value = "Tek"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1275283333333333
$ws.Range("H3").Value = 0.382585
$ws.Range("I3").Value = 0.001938448326189378
$ws.Range("J3").Value = 0.001938448326189378
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 31.083557
$ws.Range("N3").Value = 93.25067100000001
$ws.Range("O3").Value = 0.4665515560925921
$ws.Range("P3").Value = 0.4665515560925922
$ws.Range("Q3").Value = 3.964034218281667
$ws.Range("R3").Value = 35.67630796453501
$ws.Range("S3").Value = 0.0009043860829887348
$ws.Range("T3").Value = 0.000904386082988735

# Row 4: ECs -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Angptl1"
$ws.Range("C4").Value = "Tek"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1275283333333333
$ws.Range("H4").Value = 0.382585
$ws.Range("I4").Value = 0.001938448326189378
$ws.Range("J4").Value = 0.001938448326189378
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.436030666666667
$ws.Range("N4").Value = 4.308092
$ws.Range("O4").Value = 0.02155423660586901
$ws.Range("P4").Value = 0.02155423660586901
$ws.Range("Q4").Value = 0.1831345975355556
$ws.Range("R4").Value = 1.64821137782
$ws.Range("S4").Value = 0.00004178177387093659
$ws.Range("T4").Value = 0.0000417817738709366

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Angptl1"
$ws.Range("C5").Value = "Tek"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 63.34571333333334
$ws.Range("H5").Value = 190.03714
$ws.Range("I5").Value = 0.9628636144825765
$ws.Range("J5").Value = 0.9628636144825764
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 34.10446833333334
$ws.Range("N5").Value = 102.313405
$ws.Range("O5").Value = 0.5118942073015388
$ws.Range("P5").Value = 0.5118942073015389
$ws.Range("Q5").Value = 2160.371874429078
$ws.Range("R5").Value = 19443.3468698617
$ws.Range("S5").Value = 0.492884306675053
$ws.Range("T5").Value = 0.492884306675053

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Angptl1"
$ws.Range("C6").Value = "Tek"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 63.34571333333334
$ws.Range("H6").Value = 190.03714
$ws.Range("I6").Value = 0.9628636144825765
$ws.Range("J6").Value = 0.9628636144825764
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 31.083557
$ws.Range("N6").Value = 93.25067100000001
$ws.Range("O6").Value = 0.4665515560925921
$ws.Range("P6").Value = 0.4665515560925922
$ws.Range("Q6").Value = 1969.010091102327
$ws.Range("R6").Value = 17721.09081992094
$ws.Range("S6").Value = 0.4492255176417838
$ws.Range("T6").Value = 0.4492255176417839

# Row 7: FAPs -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Angptl1"
$ws.Range("C7").Value = "Tek"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 63.34571333333334
$ws.Range("H7").Value = 190.03714
$ws.Range("I7").Value = 0.9628636144825765
$ws.Range("J7").Value = 0.9628636144825764
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.436030666666667
$ws.Range("N7").Value = 4.308092
$ws.Range("O7").Value = 0.02155423660586901
$ws.Range("P7").Value = 0.02155423660586901
$ws.Range("Q7").Value = 90.96638694854224
$ws.Range("R7").Value = 818.6974825368801
$ws.Range("S7").Value = 0.02075379016573969
$ws.Range("T7").Value = 0.02075379016573969

# Row 8: sCs -> ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Angptl1"
$ws.Range("C8").Value = "Tek"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.315632666666666
$ws.Range("H8").Value = 6.946897999999999
$ws.Range("I8").Value = 0.03519793719123419
$ws.Range("J8").Value = 0.03519793719123419
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 34.10446833333334
$ws.Range("N8").Value = 102.313405
$ws.Range("O8").Value = 0.5118942073015388
$ws.Range("P8").Value = 0.5118942073015389
$ws.Range("Q8").Value = 78.97342095196555
$ws.Range("R8").Value = 710.7607885676899
$ws.Range("S8").Value = 0.01801762015715618
$ws.Range("T8").Value = 0.01801762015715618

# Row 9: sCs -> FAPs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Angptl1"
$ws.Range("C9").Value = "Tek"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.315632666666666
$ws.Range("H9").Value = 6.946897999999999
$ws.Range("I9").Value = 0.03519793719123419
$ws.Range("J9").Value = 0.03519793719123419
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 31.083557
$ws.Range("N9").Value = 93.25067100000001
$ws.Range("O9").Value = 0.4665515560925921
$ws.Range("P9").Value = 0.4665515560925922
$ws.Range("Q9").Value = 71.97809998539533
$ws.Range("R9").Value = 647.802899868558
$ws.Range("S9").Value = 0.01642165236781963
$ws.Range("T9").Value = 0.01642165236781964

# Row 10: sCs -> sCs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Angptl1"
$ws.Range("C10").Value = "Tek"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2.315632666666666
$ws.Range("H10").Value = 6.946897999999999
$ws.Range("I10").Value = 0.03519793719123419
$ws.Range("J10").Value = 0.03519793719123419
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.436030666666667
$ws.Range("N10").Value = 4.308092
$ws.Range("O10").Value = 0.02155423660586901
$ws.Range("P10").Value = 0.02155423660586901
$ws.Range("Q10").Value = 3.325319522068444
$ws.Range("R10").Value = 29.927875698616
$ws.Range("S10").Value = 0.0007586646662583782
$ws.Range("T10").Value = 0.0007586646662583783

